# Adding the "free prog users" row to the credentials sheet, matching
# the "adding free program case" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: "free prog users" -> the two extra test-account emails ---
$ws.Range("A7").Value = "free prog users"
$ws.Range("B7").Value = "rishabh.singh+1@snackmagic.com, rishabh.singh+2@snackmagic.com"

# The new row holds a long comma-separated value, so it was given wrap
# text (and, since it now wraps to two lines, a taller row height).
$ws.Range("A1:B7").WrapText = $true
$ws.Rows("7:7").RowHeight = 28.8

# Leftover selection from editing, unrelated to the data itself.
$ws.Range("G4").Select() | Out-Null
